$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.135.43'
$ws.Range("E2").Value = '  -2.18%  '
$ws.Range("D3").Value = '1.852.93'
$ws.Range("E3").Value = '  -1.15%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.41'
$ws.Range("E5").Value = '  -0.72%  '
$ws.Range("B6").Value = 'XRP'
$ws.Range("C6").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6945'
$ws.Range("E6").Value = '  -4.82%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("B8").Value = 'Dogecoin'
$ws.Range("C8").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07668'
$ws.Range("E8").Value = '  +8.06%  '
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3048'
$ws.Range("E9").Value = '  -2.64%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.45'
$ws.Range("E10").Value = '  -3.91%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08139'
$ws.Range("E11").Value = '  -1.39%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7273'
$ws.Range("E12").Value = '  -2.65%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.843.58'
$ws.Range("E13").Value = '  -3.34%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.244'
$ws.Range("E14").Value = '  -1.46%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.35'
$ws.Range("E15").Value = '  -3.30%  '
$ws.Range("D16").Value = '29.161.35'
$ws.Range("E16").Value = '  -2.13%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.795'
$ws.Range("E17").Value = '  -3.92%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007764'
$ws.Range("E18").Value = '  -0.66%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.20'
$ws.Range("E19").Value = '  -1.42%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '237.33'
$ws.Range("E20").Value = '  -4.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9998'
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("D22").Value = '2.098.18'
$ws.Range("E22").Value = '  -1.89%  '
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.626'
$ws.Range("E24").Value = '  -1.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.023'
$ws.Range("E25").Value = '  -1.61%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '161.16'
$ws.Range("E26").Value = '  -1.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1453'
$ws.Range("E27").Value = '  -5.47%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.12'
$ws.Range("E28").Value = '  -2.37%  '
$ws.Range("E29").Value = '  -2.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.400'
$ws.Range("E30").Value = '  -2.50%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.497'
$ws.Range("E31").Value = '  -0.83%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.495'
$ws.Range("E32").Value = '  -2.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.014'
$ws.Range("E33").Value = '  -4.49%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05223'
$ws.Range("E34").Value = '  -0.97%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.191'
$ws.Range("E35").Value = '  -3.59%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.033'
$ws.Range("E36").Value = '  +3.22%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.7065'
$ws.Range("E37").Value = '  -6.67%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.659'
$ws.Range("E38").Value = '  -1.71%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01865'
$ws.Range("E39").Value = '  -3.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.694'
$ws.Range("E40").Value = '  -1.96%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9365'
$ws.Range("E41").Value = '  +8.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.044'
$ws.Range("E42").Value = '  +0.61%  '
$ws.Range("D43").Value = '1.076.14'
$ws.Range("E43").Value = '  +0.72%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4287'
$ws.Range("E44").Value = '  -4.40%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '70.52'
$ws.Range("E45").Value = '  -1.25%  '
$ws.Range("E46").Value = '  -0.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '103.26'
$ws.Range("E47").Value = '  -1.36%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.781'
$ws.Range("E48").Value = '  -2.53%  '
$ws.Range("D49").Value = '1.996.85'
$ws.Range("E49").Value = '  -1.94%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.215'
$ws.Range("E50").Value = '  -3.01%  '
$ws.Range("E51").Value = '  -6.01%  '
